# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" sheet (copied from "2022-Q2" to preserve layout
#    and formatting) positioned right before "2022-Q2".
# 2. Overwrite the new sheet's fund-holding figures with the 2022-Q3 numbers.
# 3. Update the "总计" (summary) sheet: rename/retarget the first data row
#    to "2022-Q3", insert a row for the (still-present) "2022-Q2" figures,
#    and bump the running index on the "2021-Q3" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2022-Q2" -> new sheet placed immediately before it,
# then rename the duplicate to "2022-Q3".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 2: update the 2022-Q3 fund table with the new quarter's figures.
# Fund-code / name columns (B, C) are left untouched where they repeat;
# numeric-looking text columns (D-G) keep their original text formatting
# via NumberFormat "@" so leading zeros / trailing zeros survive.
# ---------------------------------------------------------------------

# Row 2: 001543 宝盈新锐灵活配置混合A
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "2.20"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "91.14"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "4.27"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0939"
$q3.Range("H2").Value = 9

# Row 3: 004448 博时汇智回报灵活配置混合 (was 007578 宝盈新锐灵活配置混合C)
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "004448"
$q3.Range("C3").Value = "博时汇智回报灵活配置混合"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "1.77"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "67.69"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "4.34"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0768"
$q3.Range("H3").Value = 5

# Row 4: 011927 博时汇誉回报混合A (was 740001 长安宏观策略混合)
$q3.Range("B4").NumberFormat = "@"
$q3.Range("B4").Value = "011927"
$q3.Range("C4").Value = "博时汇誉回报混合A"
$q3.Range("D4").NumberFormat = "@"
$q3.Range("D4").Value = "1.04"
$q3.Range("E4").NumberFormat = "@"
$q3.Range("E4").Value = "80.45"
$q3.Range("F4").NumberFormat = "@"
$q3.Range("F4").Value = "4.37"
$q3.Range("G4").NumberFormat = "@"
$q3.Range("G4").Value = "0.0454"
$q3.Range("H4").Value = 10

# Row 5: 007578 宝盈新锐灵活配置混合C (was 970073 东证融汇成长优选混合A)
$q3.Range("B5").NumberFormat = "@"
$q3.Range("B5").Value = "007578"
$q3.Range("C5").Value = "宝盈新锐灵活配置混合C"
$q3.Range("D5").NumberFormat = "@"
$q3.Range("D5").Value = "0.36"
$q3.Range("E5").NumberFormat = "@"
$q3.Range("E5").Value = "91.14"
$q3.Range("F5").NumberFormat = "@"
$q3.Range("F5").Value = "4.27"
$q3.Range("G5").NumberFormat = "@"
$q3.Range("G5").Value = "0.0154"
$q3.Range("H5").Value = 9

# Row 6: 011928 博时汇誉回报混合C (was 970074 东证融汇成长优选混合C)
$q3.Range("B6").NumberFormat = "@"
$q3.Range("B6").Value = "011928"
$q3.Range("C6").Value = "博时汇誉回报混合C"
$q3.Range("D6").NumberFormat = "@"
$q3.Range("D6").Value = "0.12"
$q3.Range("E6").NumberFormat = "@"
$q3.Range("E6").Value = "80.45"
$q3.Range("F6").NumberFormat = "@"
$q3.Range("F6").Value = "4.37"
$q3.Range("G6").NumberFormat = "@"
$q3.Range("G6").Value = "0.0052"
$q3.Range("H6").Value = 10

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Make room for the new quarter row right after the header-data row 2.
$tot.Rows.Item(3).Insert()

# Copy row 2's formatting onto the newly inserted row 3 (matches the bold
# bordered style used on column A's index cells).
$tot.Range("A2").Copy()
$tot.Range("A3").PasteSpecial(-4122)

# Row 2 now reports the new 2022-Q3 totals.
$tot.Range("B2").Value = "2022-Q3"
$tot.Range("C2").Value = 5
$tot.Range("D2").Value = 0.24

# Row 3 carries over what used to be row 2's 2022-Q2 totals.
$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2022-Q2"
$tot.Range("C3").Value = 5
$tot.Range("D3").Value = 0.17

# Row 4 (previously row 3) keeps its 2021-Q3 totals, just bumps its index.
$tot.Range("A4").Value = 2
